$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are plain numeric-looking strings (e.g. "142.88",
# "2.413.32" thousand-grouped, etc.) that must stay literal text, exactly as
# authored upstream (inline string cells). Force text format before assigning
# so Excel does not auto-convert them to numbers (which would also silently
# truncate values like "555.50" -> 555.5 or "20.30" -> 20.3).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.042.58'
$ws.Range("E2").Value = '  +2.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.416.06'
$ws.Range("E3").Value = '  +2.77%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.50'
$ws.Range("E5").Value = '  +2.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.88'
$ws.Range("E6").Value = '  +4.41%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +2.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.413.32'
$ws.Range("E9").Value = '  +2.73%  '
$ws.Range("E10").Value = '  +3.66%  '
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.18'
$ws.Range("E14").Value = '  +6.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000174'
$ws.Range("E15").Value = '  +8.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.858.95'
$ws.Range("E16").Value = '  +3.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.058.72'
$ws.Range("E17").Value = '  +2.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.417.73'
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.06'
$ws.Range("E19").Value = '  +3.79%  '
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.75'
$ws.Range("E21").Value = '  +0.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.71'
$ws.Range("E22").Value = '  +2.24%  '
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("E24").Value = '  +4.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.90'
$ws.Range("E25").Value = '  +2.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.11'
$ws.Range("E26").Value = '  +7.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '580.43'
$ws.Range("E27").Value = '  +16.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.536.64'
$ws.Range("E28").Value = '  +2.88%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +3.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0930'
$ws.Range("E31").Value = '  +7.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.46'
$ws.Range("E32").Value = '  +6.18%  '
$ws.Range("E33").Value = '  +1.02%  '
$ws.Range("E34").Value = '  +3.52%  '
$ws.Range("E35").Value = '  +2.58%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.66'
$ws.Range("E37").Value = '  +7.91%  '
$ws.Range("E38").Value = '  +3.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.383'
$ws.Range("E39").Value = '  +1.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.75'
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("E41").Value = '  +2.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '148.54'
$ws.Range("E42").Value = '  +2.47%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.75'
$ws.Range("E44").Value = '  +2.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '151.29'
$ws.Range("E45").Value = '  +6.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.29'
$ws.Range("E46").Value = '  +12.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.64'
$ws.Range("E47").Value = '  +1.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0545'
$ws.Range("E48").Value = '  +5.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.30'
$ws.Range("E49").Value = '  +6.22%  '
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("E51").Value = '  +1.62%  '
